$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data for year 2021 under the existing 2018/2019/2020 rows.
$ws.Range("A5").Value = "2021年"

# Match the formatting of the other year cells (bold, bordered, centered).
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B5").Value = -1.5
$ws.Range("C5").Value = 8.9
$ws.Range("D5").Value = -3
$ws.Range("E5").Value = 10.1
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = -2.6
